$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 946.36365
$ws.Range("I19").Value = 465.8
$ws.Range("J19").Value = 1346.8334
$ws.Range("K19").Value = 465.8
$ws.Range("L19").Value = 1346.8334
$ws.Range("M19").Value = -290.8
$ws.Range("N19").Value = -1696.8334
$ws.Range("H28").Value = 4449076.5
$ws.Range("I28").Value = 5853342
$ws.Range("K28").Value = 5853342
$ws.Range("M28").Value = -5852857
$ws.Range("H32").Value = 1199.5555
$ws.Range("I32").Value = 900
$ws.Range("J32").Value = 1237
$ws.Range("K32").Value = 900
$ws.Range("L32").Value = 1237
$ws.Range("M32").Value = -574
$ws.Range("N32").Value = -1889
$ws.Range("H40").Value = 1383.6666
$ws.Range("I40").Value = 1367.3334
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 1367.3334
$ws.Range("L40").Value = 1400
$ws.Range("M40").Value = -1192.3334
$ws.Range("N40").Value = -1750
$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3238
$ws.Range("H60").Value = 1000
$ws.Range("J60").Value = 1000
$ws.Range("L60").Value = 3000
$ws.Range("N60").Value = -3968
$ws.Range("H64").Value = 3309
$ws.Range("I64").Value = 3112.2856
$ws.Range("J64").Value = 3505.7144
$ws.Range("K64").Value = 3112.2856
$ws.Range("L64").Value = 3505.7144
$ws.Range("M64").Value = -2864.2856
$ws.Range("N64").Value = -4001.7144
$ws.Range("H67").Value = 3309
$ws.Range("I67").Value = 3112.2856
$ws.Range("J67").Value = 3505.7144
$ws.Range("K67").Value = 3112.2856
$ws.Range("L67").Value = 3505.7144
$ws.Range("M67").Value = -2254.2856
$ws.Range("N67").Value = -5221.7144
$ws.Range("H141").Value = 835.62067
$ws.Range("I141").Value = 643.0417
$ws.Range("J141").Value = 1760
$ws.Range("K141").Value = 1929.1251
$ws.Range("L141").Value = 5280
$ws.Range("M141").Value = 3250.8749
$ws.Range("N141").Value = -15640

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 703.05884
$ws.Range("I74").Value = 557.8461
$ws.Range("K74").Value = 557.8461
$ws.Range("M74").Value = 316.1539
$ws.Range("H77").Value = 703.05884
$ws.Range("I77").Value = 557.8461
$ws.Range("K77").Value = 2789.2305
$ws.Range("M77").Value = 1578.7695
$ws.Range("H97").Value = 1066.8572
$ws.Range("I97").Value = 842.96
$ws.Range("K97").Value = 842.96
$ws.Range("M97").Value = -346.96

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1204.4348
$ws.Range("I99").Value = 888.125
$ws.Range("J99").Value = 1927.4286
$ws.Range("K99").Value = 888.125
$ws.Range("L99").Value = 1927.4286
$ws.Range("M99").Value = 609.875
$ws.Range("N99").Value = -4923.4286
$ws.Range("H105").Value = 1655.7778
$ws.Range("I105").Value = 1643.2858
$ws.Range("J105").Value = 1673.2667
$ws.Range("K105").Value = 1643.2858
$ws.Range("L105").Value = 1673.2667
$ws.Range("M105").Value = 103.7141999999999
$ws.Range("N105").Value = -5167.2667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 607
$ws.Range("H70").Value = 17560
$ws.Range("I70").Value = 10080
$ws.Range("J70").Value = 18000
$ws.Range("K70").Value = 10080
$ws.Range("L70").Value = 18000
$ws.Range("M70").Value = -9765
$ws.Range("N70").Value = -18630
$ws.Range("H73").Value = 17560
$ws.Range("I73").Value = 10080
$ws.Range("J73").Value = 18000
$ws.Range("K73").Value = 10080
$ws.Range("L73").Value = 18000
$ws.Range("M73").Value = -8988
$ws.Range("N73").Value = -20184
$ws.Range("H80").Value = 21450
$ws.Range("J80").Value = 22900
$ws.Range("L80").Value = 22900
$ws.Range("N80").Value = -25146
$ws.Range("H83").Value = 21450
$ws.Range("J83").Value = 22900
$ws.Range("L83").Value = 68700
$ws.Range("N83").Value = -79932
$ws.Range("H86").Value = 33111.38
$ws.Range("I86").Value = 6860.9443
$ws.Range("J86").Value = 76066.63
$ws.Range("K86").Value = 6860.9443
$ws.Range("L86").Value = 76066.63
$ws.Range("M86").Value = -5737.9443
$ws.Range("N86").Value = -78312.63
$ws.Range("H89").Value = 33111.38
$ws.Range("I89").Value = 6860.9443
$ws.Range("J89").Value = 76066.63
$ws.Range("K89").Value = 34304.7215
$ws.Range("L89").Value = 380333.15
$ws.Range("M89").Value = -28688.7215
$ws.Range("N89").Value = -391565.15
$ws.Range("H104").Value = 23326.666
$ws.Range("J104").Value = 23326.666
$ws.Range("L104").Value = 23326.666
$ws.Range("N104").Value = -28568.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 415
$ws.Range("J50").Value = 580.9091
$ws.Range("L50").Value = 1742.7273
$ws.Range("N50").Value = -2704.7273
$ws.Range("H53").Value = 415
$ws.Range("J53").Value = 580.9091
$ws.Range("L53").Value = 1742.7273
$ws.Range("N53").Value = -2704.7273
$ws.Range("H64").Value = 540.6667
$ws.Range("I64").Value = 248.8
$ws.Range("K64").Value = 746.4000000000001
$ws.Range("M64").Value = -476.4000000000001
$ws.Range("H67").Value = 540.6667
$ws.Range("I67").Value = 248.8
$ws.Range("K67").Value = 746.4000000000001
$ws.Range("M67").Value = 189.5999999999999
$ws.Range("H88").Value = 5200
$ws.Range("J88").Value = 5200
$ws.Range("L88").Value = 15600
$ws.Range("N88").Value = -16456
$ws.Range("H91").Value = 5200
$ws.Range("J91").Value = 5200
$ws.Range("L91").Value = 15600
$ws.Range("N91").Value = -18564

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7938303.5
$ws.Range("I126").Value = 1692.4615
$ws.Range("J126").Value = 20835296
$ws.Range("K126").Value = 5077.3845
$ws.Range("L126").Value = 62505888
$ws.Range("M126").Value = -2607.3845
$ws.Range("N126").Value = -62510828

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2449.9473
$ws.Range("I7").Value = 1682.0714
$ws.Range("J7").Value = 4600
$ws.Range("K7").Value = 1682.0714
$ws.Range("L7").Value = 4600
$ws.Range("M7").Value = -1570.0714
$ws.Range("N7").Value = -4824
$ws.Range("H40").Value = 1742.1786
$ws.Range("I40").Value = 1608.381
$ws.Range("J40").Value = 2143.5715
$ws.Range("K40").Value = 1608.381
$ws.Range("L40").Value = 2143.5715
$ws.Range("M40").Value = -1472.381
$ws.Range("N40").Value = -2415.5715
$ws.Range("H47").Value = 5399
$ws.Range("J47").Value = 5399
$ws.Range("L47").Value = 5399
$ws.Range("N47").Value = -6379
$ws.Range("H52").Value = 5399
$ws.Range("J52").Value = 5399
$ws.Range("L52").Value = 5399
$ws.Range("N52").Value = -5865
$ws.Range("H126").Value = 2449.9473
$ws.Range("I126").Value = 1682.0714
$ws.Range("J126").Value = 4600
$ws.Range("K126").Value = 5046.2142
$ws.Range("L126").Value = 13800
$ws.Range("M126").Value = -2576.2142
$ws.Range("N126").Value = -18740
$ws.Range("H128").Value = 41963
$ws.Range("J128").Value = 41963
$ws.Range("L128").Value = 41963
$ws.Range("N128").Value = -51923
$ws.Range("H136").Value = 4558.135
$ws.Range("I136").Value = 4206.2104
$ws.Range("K136").Value = 12618.6312
$ws.Range("M136").Value = -10068.6312

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 530.5909
$ws.Range("J107").Value = 588.6667
$ws.Range("L107").Value = 1766.0001
$ws.Range("N107").Value = -5606.0001
